$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.345.76"
$ws.Range("E2").Value = "  -3.39%  "

$ws.Range("D3").Value = "2.476.41"
$ws.Range("E3").Value = "  -2.58%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "312.71"
$ws.Range("E5").Value = "  +0.24%  "

$ws.Range("D6").Value = "94.32"
$ws.Range("E6").Value = "  -5.91%  "

$ws.Range("D7").Value = "0.549"
$ws.Range("E7").Value = "  -3.17%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -4.73%  "

$ws.Range("E10").Value = "  -5.59%  "

$ws.Range("D11").Value = "0.0782"
$ws.Range("E11").Value = "  -2.91%  "

$ws.Range("E12").Value = "  -1.00%  "

$ws.Range("D13").Value = "6.99"
$ws.Range("E13").Value = "  -4.35%  "

$ws.Range("D14").Value = "2.855.88"
$ws.Range("E14").Value = "  -2.69%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.490.56"
$ws.Range("E15").Value = "  -4.88%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "15.25"
$ws.Range("E16").Value = "  -1.07%  "

$ws.Range("E17").Value = "  -3.33%  "

$ws.Range("D18").Value = "41.322.07"
$ws.Range("E18").Value = "  -3.45%  "

$ws.Range("E19").Value = "  -6.30%  "

$ws.Range("D20").Value = "0.0₃0926"
$ws.Range("E20").Value = "  -2.64%  "

$ws.Range("D21").Value = "11.21"
$ws.Range("E21").Value = "  -9.30%  "

$ws.Range("D22").Value = "68.68"
$ws.Range("E22").Value = "  -1.73%  "

$ws.Range("D23").Value = "236.99"
$ws.Range("E23").Value = "  -2.49%  "

$ws.Range("E24").Value = "  -4.74%  "

$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("E26").Value = "  -6.14%  "

$ws.Range("D27").Value = "24.06"
$ws.Range("E27").Value = "  -6.35%  "

$ws.Range("E28").Value = "  -4.44%  "

$ws.Range("E29").Value = "  -4.46%  "

$ws.Range("E30").Value = "  -5.37%  "

$ws.Range("D31").Value = "151.71"
$ws.Range("E31").Value = "  -4.06%  "

$ws.Range("E32").Value = "  -7.72%  "

$ws.Range("E33").Value = "  -3.83%  "

$ws.Range("D34").Value = "2.55"
$ws.Range("E34").Value = "  -7.04%  "

$ws.Range("E35").Value = "  -6.27%  "

$ws.Range("E36").Value = "  -2.94%  "

$ws.Range("D37").Value = "17.41"
$ws.Range("E37").Value = "  -2.71%  "

$ws.Range("E38").Value = "  -5.19%  "

$ws.Range("E39").Value = "  -2.58%  "

$ws.Range("E40").Value = "  -8.93%  "

$ws.Range("E41").Value = "  +1.41%  "

$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.25%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "19.47"
$ws.Range("E43").Value = "  -10.95%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.990.49"
$ws.Range("E44").Value = "  -0.36%  "

$ws.Range("E45").Value = "  -4.68%  "

$ws.Range("E46").Value = "  -8.36%  "

$ws.Range("D47").Value = "8.70"
$ws.Range("E47").Value = "  -4.16%  "

$ws.Range("D48").Value = "2.719.95"
$ws.Range("E48").Value = "  -2.45%  "

$ws.Range("D49").Value = "69.43"
$ws.Range("E49").Value = "  -4.01%  "

$ws.Range("D50").Value = "97.11"
$ws.Range("E50").Value = "  -4.33%  "

$ws.Range("D51").Value = "74.38"
$ws.Range("E51").Value = "  -6.82%  "
